{"js": "// This document contains a single table of two-digit \u00f7 one-digit\n// division \"problem, quotient, remainder\" strings. Data only lives in\n// every 4th row (0-based rows 0, 4, 8, 12, 16); the rows in between are\n// blank spacer rows, and every data row has 5 populated columns (0-4).\n//\n// Each cell is targeted *positionally* (table row/col) rather than via\n// a document-wide text search: one cell's new value happens to equal\n// another cell's old value (\"81\u00f75=16, 1\"), so a global find-and-replace\n// executed in document order would end up matching \u2014 and corrupting \u2014\n// the wrong cell once the first one had already been rewritten.\n// Scoping the search to the individual cell's own body sidesteps that\n// ambiguity entirely, and using search()+insertText(..., \"Replace\") on\n// the matched range (instead of clearing/resetting the cell body)\n// preserves the existing run/paragraph formatting (fonts, size,\n// alignment).\n\nconst replacements = [\n  { row: 0, col: 0, old: \"50\u00f79=5, 5\", new: \"79\u00f77=11, 2\" },\n  { row: 0, col: 1, old: \"49\u00f72=24, 1\", new: \"19\u00f76=3, 1\" },\n  { row: 0, col: 2, old: \"90\u00f79=10, 0\", new: \"10\u00f75=2, 0\" },\n  { row: 0, col: 3, old: \"49\u00f77=7, 0\", new: \"46\u00f74=11, 2\" },\n  { row: 0, col: 4, old: \"23\u00f77=3, 2\", new: \"97\u00f75=19, 2\" },\n  { row: 4, col: 0, old: \"73\u00f79=8, 1\", new: \"96\u00f76=16, 0\" },\n  { row: 4, col: 1, old: \"52\u00f74=13, 0\", new: \"39\u00f77=5, 4\" },\n  { row: 4, col: 2, old: \"47\u00f75=9, 2\", new: \"89\u00f77=12, 5\" },\n  { row: 4, col: 3, old: \"71\u00f75=14, 1\", new: \"97\u00f72=48, 1\" },\n  { row: 4, col: 4, old: \"60\u00f78=7, 4\", new: \"41\u00f74=10, 1\" },\n  { row: 8, col: 0, old: \"11\u00f72=5, 1\", new: \"91\u00f78=11, 3\" },\n  { row: 8, col: 1, old: \"84\u00f77=12, 0\", new: \"81\u00f75=16, 1\" },\n  { row: 8, col: 2, old: \"46\u00f79=5, 1\", new: \"21\u00f76=3, 3\" },\n  { row: 8, col: 3, old: \"31\u00f74=7, 3\", new: \"18\u00f73=6, 0\" },\n  { row: 8, col: 4, old: \"78\u00f79=8, 6\", new: \"78\u00f76=13, 0\" },\n  { row: 12, col: 0, old: \"74\u00f73=24, 2\", new: \"38\u00f72=19, 0\" },\n  { row: 12, col: 1, old: \"35\u00f73=11, 2\", new: \"48\u00f79=5, 3\" },\n  { row: 12, col: 2, old: \"40\u00f76=6, 4\", new: \"57\u00f79=6, 3\" },\n  { row: 12, col: 3, old: \"10\u00f76=1, 4\", new: \"72\u00f78=9, 0\" },\n  { row: 12, col: 4, old: \"42\u00f75=8, 2\", new: \"90\u00f78=11, 2\" },\n  { row: 16, col: 0, old: \"17\u00f78=2, 1\", new: \"78\u00f78=9, 6\" },\n  { row: 16, col: 1, old: \"58\u00f79=6, 4\", new: \"60\u00f75=12, 0\" },\n  { row: 16, col: 2, old: \"20\u00f76=3, 2\", new: \"11\u00f76=1, 5\" },\n  { row: 16, col: 3, old: \"81\u00f75=16, 1\", new: \"92\u00f79=10, 2\" },\n  { row: 16, col: 4, old: \"74\u00f76=12, 2\", new: \"20\u00f78=2, 4\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body\");\n}\nconst table = tables.items[0];\n\nfor (const rep of replacements) {\n  const cell = table.getCell(rep.row, rep.col);\n  const results = cell.body.search(rep.old, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${rep.old}\" in cell (${rep.row},${rep.col}), found ${results.items.length}`\n    );\n  }\n  results.items[0].insertText(rep.new, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# This document contains a single table of two-digit / one-digit\n# division \"problem, quotient, remainder\" strings. Data only lives in\n# every 4th row (Word's 1-based rows 1, 5, 9, 13, 17); the rows in\n# between are blank spacer rows, and every data row has 5 populated\n# columns (1-5).\n#\n# Each cell is targeted *positionally* via Table.Cell(row, col) rather\n# than a document-wide Find/Replace: one cell's new value happens to\n# equal another cell's old value (\"81\u00f75=16, 1\"), so a global replace-all\n# executed in document order would end up matching -- and corrupting --\n# the wrong cell once the first one had already been rewritten.\n# Addressing cells by position sidesteps that ambiguity entirely. Each\n# cell's current text is verified against the expected old value before\n# being overwritten (defensive check), and only the cell Range's text is\n# reassigned, which preserves the run/paragraph formatting already on\n# that range (fonts, size, alignment).\n\n$d = $word.ActiveDocument\nif ($d.Tables.Count -eq 0) {\n    throw \"No table found in document\"\n}\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1; Col = 1; Old = \"50\u00f79=5, 5\"; New = \"79\u00f77=11, 2\" }\n    @{ Row = 1; Col = 2; Old = \"49\u00f72=24, 1\"; New = \"19\u00f76=3, 1\" }\n    @{ Row = 1; Col = 3; Old = \"90\u00f79=10, 0\"; New = \"10\u00f75=2, 0\" }\n    @{ Row = 1; Col = 4; Old = \"49\u00f77=7, 0\"; New = \"46\u00f74=11, 2\" }\n    @{ Row = 1; Col = 5; Old = \"23\u00f77=3, 2\"; New = \"97\u00f75=19, 2\" }\n    @{ Row = 5; Col = 1; Old = \"73\u00f79=8, 1\"; New = \"96\u00f76=16, 0\" }\n    @{ Row = 5; Col = 2; Old = \"52\u00f74=13, 0\"; New = \"39\u00f77=5, 4\" }\n    @{ Row = 5; Col = 3; Old = \"47\u00f75=9, 2\"; New = \"89\u00f77=12, 5\" }\n    @{ Row = 5; Col = 4; Old = \"71\u00f75=14, 1\"; New = \"97\u00f72=48, 1\" }\n    @{ Row = 5; Col = 5; Old = \"60\u00f78=7, 4\"; New = \"41\u00f74=10, 1\" }\n    @{ Row = 9; Col = 1; Old = \"11\u00f72=5, 1\"; New = \"91\u00f78=11, 3\" }\n    @{ Row = 9; Col = 2; Old = \"84\u00f77=12, 0\"; New = \"81\u00f75=16, 1\" }\n    @{ Row = 9; Col = 3; Old = \"46\u00f79=5, 1\"; New = \"21\u00f76=3, 3\" }\n    @{ Row = 9; Col = 4; Old = \"31\u00f74=7, 3\"; New = \"18\u00f73=6, 0\" }\n    @{ Row = 9; Col = 5; Old = \"78\u00f79=8, 6\"; New = \"78\u00f76=13, 0\" }\n    @{ Row = 13; Col = 1; Old = \"74\u00f73=24, 2\"; New = \"38\u00f72=19, 0\" }\n    @{ Row = 13; Col = 2; Old = \"35\u00f73=11, 2\"; New = \"48\u00f79=5, 3\" }\n    @{ Row = 13; Col = 3; Old = \"40\u00f76=6, 4\"; New = \"57\u00f79=6, 3\" }\n    @{ Row = 13; Col = 4; Old = \"10\u00f76=1, 4\"; New = \"72\u00f78=9, 0\" }\n    @{ Row = 13; Col = 5; Old = \"42\u00f75=8, 2\"; New = \"90\u00f78=11, 2\" }\n    @{ Row = 17; Col = 1; Old = \"17\u00f78=2, 1\"; New = \"78\u00f78=9, 6\" }\n    @{ Row = 17; Col = 2; Old = \"58\u00f79=6, 4\"; New = \"60\u00f75=12, 0\" }\n    @{ Row = 17; Col = 3; Old = \"20\u00f76=3, 2\"; New = \"11\u00f76=1, 5\" }\n    @{ Row = 17; Col = 4; Old = \"81\u00f75=16, 1\"; New = \"92\u00f79=10, 2\" }\n    @{ Row = 17; Col = 5; Old = \"74\u00f76=12, 2\"; New = \"20\u00f78=2, 4\" }\n)\n\nforeach ($rep in $replacements) {\n    $cell = $t.Cell($rep.Row, $rep.Col)\n    $current = $cell.Range.Text.TrimEnd([char]7, [char]13)\n    if ($current -ne $rep.Old) {\n        throw \"Cell ($($rep.Row),$($rep.Col)) expected `\"$($rep.Old)`\" but found `\"$current`\"\"\n    }\n    $cell.Range.Text = $rep.New\n}\n"}
